# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Estado de Cuenta" detail table (rows 16-43, columns C:G) is
# re-sorted/refreshed: each worker's overdue-period rows are regrouped
# together and ordered by period (most recent period first), with the
# matching "Valor Mora" (F) / "Salario Basico" (G) values carried along
# with their (document, period) pair. Rows 16-28 become
# "45523211 / KELLY NUÑEZ AVILA" (periods 2102 down to 2002) and rows
# 29-43 become "1047388104 / ANA GABRIEL GAVALO EMITOLA" (periods 2102
# down to 1912).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2102"; F = 25749; G = 877803 },
    @{ Row = 17; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2101"; F = 35112; G = 877803 },
    @{ Row = 18; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2012"; F = 35112; G = 877803 },
    @{ Row = 19; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2011"; F = 35112; G = 877803 },
    @{ Row = 20; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2010"; F = 35112; G = 877803 },
    @{ Row = 21; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2009"; F = 35112; G = 877803 },
    @{ Row = 22; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2008"; F = 35112; G = 877803 },
    @{ Row = 23; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2007"; F = 35112; G = 877803 },
    @{ Row = 24; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2006"; F = 35112; G = 877803 },
    @{ Row = 25; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2005"; F = 35112; G = 877803 },
    @{ Row = 26; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2004"; F = 35112; G = 877803 },
    @{ Row = 27; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2003"; F = 35112; G = 877803 },
    @{ Row = 28; C = "45523211"; D = "KELLY NUÑEZ AVILA"; E = "2002"; F = 35112; G = 877803 },
    @{ Row = 29; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2102"; F = 24292; G = 828116 },
    @{ Row = 30; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2101"; F = 33125; G = 828116 },
    @{ Row = 31; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2012"; F = 33125; G = 828116 },
    @{ Row = 32; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2011"; F = 33125; G = 828116 },
    @{ Row = 33; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2010"; F = 33125; G = 828116 },
    @{ Row = 34; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2009"; F = 33125; G = 828116 },
    @{ Row = 35; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2008"; F = 33125; G = 828116 },
    @{ Row = 36; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2007"; F = 33125; G = 828116 },
    @{ Row = 37; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2006"; F = 33125; G = 828116 },
    @{ Row = 38; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2005"; F = 33125; G = 828116 },
    @{ Row = 39; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2004"; F = 33125; G = 828116 },
    @{ Row = 40; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2003"; F = 33125; G = 828116 },
    @{ Row = 41; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2002"; F = 33125; G = 828116 },
    @{ Row = 42; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "2001"; F = 33125; G = 828116 },
    @{ Row = 43; C = "1047388104"; D = "ANA GABRIEL GAVALO EMITOLA"; E = "1912"; F = 33125; G = 828116 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
}
